$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1844370.4
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1844370.4
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5533111.199999999
$ws.Range("N17").Value = -5533447.199999999
# Row 32
$ws.Range("H32").Value = 275.85715
$ws.Range("I32").Value = 282.25
$ws.Range("J32").Value = 267.33334
$ws.Range("K32").Value = 282.25
$ws.Range("L32").Value = 267.33334
$ws.Range("M32").Value = 43.75
$ws.Range("N32").Value = -919.33334
# Row 98
$ws.Range("H98").Value = 2445.9167
$ws.Range("I98").Value = 1982.7778
$ws.Range("J98").Value = 3835.3333
$ws.Range("K98").Value = 1982.7778
$ws.Range("L98").Value = 3835.3333
$ws.Range("M98").Value = -484.7778000000001
$ws.Range("N98").Value = -6831.3333
# Row 112
$ws.Range("H112").Value = 1165.6897
$ws.Range("I112").Value = 966.6667
$ws.Range("J112").Value = 1188.6538
$ws.Range("K112").Value = 2900.0001
$ws.Range("L112").Value = 3565.9614
$ws.Range("M112").Value = -1792.0001
$ws.Range("N112").Value = -5781.9614
# Row 116
$ws.Range("H116").Value = 6532.7085
$ws.Range("I116").Value = 9398.846
$ws.Range("J116").Value = 3145.4546
$ws.Range("K116").Value = 9398.846
$ws.Range("L116").Value = 3145.4546
$ws.Range("M116").Value = -5956.846
$ws.Range("N116").Value = -10029.4546
# Row 122
$ws.Range("H122").Value = 2445.9167
$ws.Range("I122").Value = 1982.7778
$ws.Range("J122").Value = 3835.3333
$ws.Range("K122").Value = 5948.3334
$ws.Range("L122").Value = 11505.9999
$ws.Range("M122").Value = -3498.3334
$ws.Range("N122").Value = -16405.9999
# Row 132
$ws.Range("H132").Value = 1443.2667
$ws.Range("I132").Value = 1062.1562
$ws.Range("J132").Value = 2381.3845
$ws.Range("K132").Value = 3186.4686
$ws.Range("L132").Value = 7144.1535
$ws.Range("M132").Value = -656.4685999999997
$ws.Range("N132").Value = -12204.1535
# Row 137
$ws.Range("H137").Value = 1416.6182
$ws.Range("I137").Value = 1177.275
$ws.Range("J137").Value = 2054.8667
$ws.Range("K137").Value = 3531.825
$ws.Range("L137").Value = 6164.6001
$ws.Range("M137").Value = -981.8250000000003
$ws.Range("N137").Value = -11264.6001
# Row 138
$ws.Range("H138").Value = 2142.282
$ws.Range("I138").Value = 810.42
$ws.Range("J138").Value = 4520.607
$ws.Range("K138").Value = 2431.26
$ws.Range("L138").Value = 13561.821
$ws.Range("M138").Value = 2708.74
$ws.Range("N138").Value = -23841.821

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1287.6471
$ws.Range("I61").Value = 1140.4103
$ws.Range("J61").Value = 1766.1666
$ws.Range("K61").Value = 1140.4103
$ws.Range("L61").Value = 1766.1666
$ws.Range("M61").Value = -928.4103
$ws.Range("N61").Value = -2190.1666
# Row 74
$ws.Range("H74").Value = 1338.8
$ws.Range("I74").Value = 1112.88
$ws.Range("J74").Value = 2468.4
$ws.Range("K74").Value = 1112.88
$ws.Range("L74").Value = 2468.4
$ws.Range("M74").Value = -238.8800000000001
$ws.Range("N74").Value = -4216.4
# Row 77
$ws.Range("H77").Value = 1338.8
$ws.Range("I77").Value = 1112.88
$ws.Range("J77").Value = 2468.4
$ws.Range("K77").Value = 5564.400000000001
$ws.Range("L77").Value = 12342
$ws.Range("M77").Value = -1196.400000000001
$ws.Range("N77").Value = -21078
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
# Row 101
$ws.Range("H101").Value = 60000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 60000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 60000
$ws.Range("N101").Value = -66490
# Row 122
$ws.Range("H122").Value = 4279373
$ws.Range("I122").Value = 6418309
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 19254927
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -19252477
$ws.Range("N122").Value = -9400
# Row 136
$ws.Range("H136").Value = 1287.6471
$ws.Range("I136").Value = 1140.4103
$ws.Range("J136").Value = 1766.1666
$ws.Range("K136").Value = 3421.2309
$ws.Range("L136").Value = 5298.4998
$ws.Range("M136").Value = -871.2309
$ws.Range("N136").Value = -10398.4998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value = 21900
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 21900
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 21900
$ws.Range("N88").Value = -22712
# Row 91
$ws.Range("H91").Value = 21900
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 21900
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 21900
$ws.Range("N91").Value = -24708
# Row 105
$ws.Range("H105").Value = 26328536
$ws.Range("I105").Value = 62526510
$ws.Range("J105").Value = 2736.182
$ws.Range("K105").Value = 62526510
$ws.Range("L105").Value = 2736.182
$ws.Range("M105").Value = -62524763
$ws.Range("N105").Value = -6230.182
# Row 126
$ws.Range("H126").Value = 62217.125
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 62217.125
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 62217.125
$ws.Range("N126").Value = -72097.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 2086.818
$ws.Range("I134").Value = 2489.5881
$ws.Range("J134").Value = 1434.7142
$ws.Range("K134").Value = 7468.7643
$ws.Range("L134").Value = 4304.142599999999
$ws.Range("M134").Value = -4933.7643
$ws.Range("N134").Value = -9374.142599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 2439899.2
$ws.Range("I131").Value = 5882959
$ws.Range("J131").Value = 1065.4584
$ws.Range("K131").Value = 17648877
$ws.Range("L131").Value = 3196.3752
$ws.Range("M131").Value = -17643837
$ws.Range("N131").Value = -13276.3752

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2233.0857
$ws.Range("I132").Value = 2005.2727
$ws.Range("J132").Value = 2618.6155
$ws.Range("K132").Value = 6015.8181
$ws.Range("L132").Value = 7855.8465
$ws.Range("M132").Value = -3485.8181
$ws.Range("N132").Value = -12915.8465

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 3737.2075
$ws.Range("I136").Value = 1669.5227
$ws.Range("J136").Value = 13845.889
$ws.Range("K136").Value = 5008.5681
$ws.Range("L136").Value = 41537.667
$ws.Range("M136").Value = -2458.5681
$ws.Range("N136").Value = -46637.667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 534
$ws.Range("I100").Value = 401
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 802
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -261
$ws.Range("N100").Value = -2682
# Row 132
$ws.Range("H132").Value = 28031.621
$ws.Range("I132").Value = 33015.324
$ws.Range("J132").Value = 2282.5
$ws.Range("K132").Value = 99045.97200000001
$ws.Range("L132").Value = 6847.5
$ws.Range("M132").Value = -96515.97200000001
$ws.Range("N132").Value = -11907.5
# Row 136
$ws.Range("H136").Value = 7938888.5
$ws.Range("I136").Value = 2635.7073
$ws.Range("J136").Value = 22729178
$ws.Range("K136").Value = 7907.1219
$ws.Range("L136").Value = 68187534
$ws.Range("M136").Value = -5357.1219
$ws.Range("N136").Value = -68192634

$wb.Save()
Write-Host "Done applying edits"